# Update the "Handback" report timestamps for the
# a35f5189-d9b8-4f18-a057-abf087767076 entry, as produced by a re-run of
# the handback status report generation.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the a35f5189... row (row 4)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-21 18:49:55"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the a35f5189... row (row 4)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-21 18:49:50"
$wsZhCn.Range("K4").Value = "2016-08-21 18:50:16"

# de-de sheet: "Correspond Handback DateTime" for the a35f5189... row (row 4)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-08-21 18:50:22"
